# Auto-generated edit script: applies the numeric corrections from the
# 'chore: update Sheets via scheduled runner' commit to each affected sheet.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 39000
$ws.Range("J3").Value = 39000
$ws.Range("L3").Value = 39000
$ws.Range("N3").Value = -39228
$ws.Range("H87").Value = 18787.82
$ws.Range("J87").Value = 18787.82
$ws.Range("L87").Value = 18787.82
$ws.Range("N87").Value = -21283.82
$ws.Range("H90").Value = 18787.82
$ws.Range("J90").Value = 18787.82
$ws.Range("L90").Value = 56363.46
$ws.Range("N90").Value = -68843.45999999999
$ws.Range("H102").Value = 39000
$ws.Range("J102").Value = 39000
$ws.Range("L102").Value = 39000
$ws.Range("N102").Value = -45490

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H30").Value = 1000
$ws.Range("I30").Value = 1000
$ws.Range("K30").Value = 1000
$ws.Range("M30").Value = -850
$ws.Range("H32").Value = 4104.414
$ws.Range("I32").Value = 3548.4111
$ws.Range("K32").Value = 3548.4111
$ws.Range("M32").Value = -3261.4111
$ws.Range("H43").Value = 15030.8
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 15030.8
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 15030.8
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = -15656.8
$ws.Range("H74").Value = 3977.5
$ws.Range("I74").Value = 4538.515
$ws.Range("J74").Value = 1332.7142
$ws.Range("K74").Value = 4538.515
$ws.Range("L74").Value = 1332.7142
$ws.Range("M74").Value = -3664.515
$ws.Range("N74").Value = -3080.7142
$ws.Range("H77").Value = 3977.5
$ws.Range("I77").Value = 4538.515
$ws.Range("J77").Value = 1332.7142
$ws.Range("K77").Value = 22692.575
$ws.Range("L77").Value = 6663.571
$ws.Range("M77").Value = -18324.575
$ws.Range("N77").Value = -15399.571
$ws.Range("H109").Value = 24710.666
$ws.Range("J109").Value = 24710.666
$ws.Range("L109").Value = 24710.666
$ws.Range("N109").Value = -27484.666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H110").Value = 31867.334
$ws.Range("J110").Value = 31867.334
$ws.Range("L110").Value = 31867.334
$ws.Range("N110").Value = -40047.334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 4951.25
$ws.Range("I132").Value = 2204
$ws.Range("J132").Value = 6599.6
$ws.Range("K132").Value = 6612
$ws.Range("L132").Value = 19798.8
$ws.Range("M132").Value = -4082
$ws.Range("N132").Value = -24858.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1600
$ws.Range("I4").Value = 400
$ws.Range("J4").Value = 2000
$ws.Range("K4").Value = 1200
$ws.Range("L4").Value = 6000
$ws.Range("M4").Value = -1088
$ws.Range("N4").Value = -6224
$ws.Range("H5").Value = 924664.3
$ws.Range("I5").Value = 761.0909
$ws.Range("J5").Value = 2195031.2
$ws.Range("K5").Value = 2283.2727
$ws.Range("L5").Value = 6585093.600000001
$ws.Range("M5").Value = -2171.2727
$ws.Range("N5").Value = -6585317.600000001
$ws.Range("H63").Value = 3350.647
$ws.Range("I63").Value = 2018.5
$ws.Range("J63").Value = 3528.2666
$ws.Range("K63").Value = 6055.5
$ws.Range("L63").Value = 10584.7998
$ws.Range("M63").Value = -5306.5
$ws.Range("N63").Value = -12082.7998
$ws.Range("H66").Value = 3350.647
$ws.Range("I66").Value = 2018.5
$ws.Range("J66").Value = 3528.2666
$ws.Range("K66").Value = 18166.5
$ws.Range("L66").Value = 31754.3994
$ws.Range("M66").Value = -14422.5
$ws.Range("N66").Value = -39242.39939999999
$ws.Range("H121").Value = 309199.12
$ws.Range("I121").Value = 166939.67
$ws.Range("J121").Value = 337651
$ws.Range("K121").Value = 500819.01
$ws.Range("L121").Value = 1012953
$ws.Range("M121").Value = -499509.01
$ws.Range("N121").Value = -1015573
$ws.Range("H129").Value = 1767.4667
$ws.Range("J129").Value = 2635.2354
$ws.Range("L129").Value = 7905.706200000001
$ws.Range("N129").Value = -17905.7062
$ws.Range("H131").Value = 4388.409
$ws.Range("I131").Value = 760
$ws.Range("J131").Value = 4961.316
$ws.Range("K131").Value = 2280
$ws.Range("L131").Value = 14883.948
$ws.Range("M131").Value = 2760
$ws.Range("N131").Value = -24963.948
$ws.Range("H135").Value = 924664.3
$ws.Range("I135").Value = 761.0909
$ws.Range("J135").Value = 2195031.2
$ws.Range("K135").Value = 6849.8181
$ws.Range("L135").Value = 19755280.8
$ws.Range("M135").Value = -4314.8181
$ws.Range("N135").Value = -19760350.8
$ws.Range("H137").Value = 2194.6296
$ws.Range("I137").Value = 2702.111
$ws.Range("J137").Value = 1940.8889
$ws.Range("K137").Value = 8106.333
$ws.Range("L137").Value = 5822.6667
$ws.Range("M137").Value = -3006.333
$ws.Range("N137").Value = -16022.6667
$ws.Range("H139").Value = 1675.7826
$ws.Range("I139").Value = 1240.2727
$ws.Range("J139").Value = 2075
$ws.Range("K139").Value = 3720.8181
$ws.Range("L139").Value = 6225
$ws.Range("M139").Value = 1419.1819
$ws.Range("N139").Value = -16505

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2495.7585
$ws.Range("I122").Value = 1989.9333
$ws.Range("K122").Value = 5969.7999
$ws.Range("M122").Value = -3519.7999
$ws.Range("H126").Value = 2069
$ws.Range("I126").Value = 1851.1333
$ws.Range("K126").Value = 5553.3999
$ws.Range("M126").Value = -3083.3999
$ws.Range("H132").Value = 1873.9
$ws.Range("I132").Value = 1498.25
$ws.Range("J132").Value = 3376.5
$ws.Range("K132").Value = 4494.75
$ws.Range("L132").Value = 10129.5
$ws.Range("M132").Value = -1964.75
$ws.Range("N132").Value = -15189.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4221.979
$ws.Range("I132").Value = 4074.419
$ws.Range("J132").Value = 4718.3184
$ws.Range("K132").Value = 12223.257
$ws.Range("L132").Value = 14154.9552
$ws.Range("M132").Value = -9693.257
$ws.Range("N132").Value = -19214.9552
$ws.Range("H136").Value = 7753380
$ws.Range("I136").Value = 1472.3429
$ws.Range("J136").Value = 41667976
$ws.Range("K136").Value = 4417.028700000001
$ws.Range("L136").Value = 125003928
$ws.Range("M136").Value = -1867.028700000001
$ws.Range("N136").Value = -125009028

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 50599.816
$ws.Range("J46").Value = 50599.816
$ws.Range("L46").Value = 50599.816
$ws.Range("N46").Value = -51061.816
$ws.Range("H122").Value = 41618.88
$ws.Range("I122").Value = 54009.844
$ws.Range("K122").Value = 162029.532
$ws.Range("M122").Value = -159579.532
$ws.Range("H134").Value = 50599.816
$ws.Range("J134").Value = 50599.816
$ws.Range("L134").Value = 151799.448
$ws.Range("N134").Value = -156869.448
$ws.Range("H136").Value = 3733.2888
$ws.Range("I136").Value = 651.4074000000001
$ws.Range("J136").Value = 8356.111000000001
$ws.Range("K136").Value = 1954.2222
$ws.Range("L136").Value = 25068.333
$ws.Range("M136").Value = 595.7777999999998
$ws.Range("N136").Value = -30168.333
